$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Find-ParaByText($doc, $needle) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------
# 1) "Stalon Vilhelminas kommun Västerbotten" -> split "Stalon" out
#    with spell-check proofErr markers.
# ---------------------------------------------------------------
$p = Find-ParaByText $d "Stalon Vilhelminas kommun"
$xml = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:pPr>' +
    '<w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr>' +
    '<w:spacing w:before="0"/>' +
    '<w:rPr><w:rFonts w:ascii="PT Sans Narrow" w:eastAsia="PT Sans Narrow" w:hAnsi="PT Sans Narrow" w:cs="PT Sans Narrow"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
    '</w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="PT Sans Narrow" w:eastAsia="PT Sans Narrow" w:hAnsi="PT Sans Narrow" w:cs="PT Sans Narrow"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Stalon</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="PT Sans Narrow" w:eastAsia="PT Sans Narrow" w:hAnsi="PT Sans Narrow" w:cs="PT Sans Narrow"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> Vilhelminas kommun Västerbotten</w:t></w:r>' +
    '</w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------
# 2) "Målet med webbplatsen ... kist priser och våra schyssta  "
#    -> split "kist" out with proofErr, drop the _GoBack bookmark
#    here, and append a new "el installationer " run.
# ---------------------------------------------------------------
$p = Find-ParaByText $d "Målet med webbplatsen"
$xml = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:r><w:t xml:space="preserve">Målet med webbplatsen är att kunna få reda på våra låga </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>kist</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> priser och våra</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> schyssta </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">el installationer </w:t></w:r>' +
    '</w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------
# 3) The empty paragraph right after "Målgrupp" gets the text
#    "24-45".
# ---------------------------------------------------------------
$p = Find-ParaByText $d "Målgrupp"
$paras = $d.Paragraphs
$emptyPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -eq "Målgrupp`r") {
        $emptyPara = $paras.Item($i + 1)
        break
    }
}
$xml = '<w:p xmlns:w="' + $wNs + '"><w:r><w:t>24-45</w:t></w:r></w:p>'
$emptyPara.Range.InsertXML($xml)

# ---------------------------------------------------------------
# 4) New paragraph after "Research och inhämtning av material":
#    "Lätt att hitta. Telefon nummer högst uppe" + ". Facebook sida . "
# ---------------------------------------------------------------
$p = Find-ParaByText $d "Research och inhämtning av material"
$p.Range.InsertParagraphAfter()
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -eq "Research och inhämtning av material`r") {
        $newPara = $paras.Item($i + 1)
        break
    }
}
$xml = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:r><w:t>Lätt att hitta. Telefon nummer högst uppe</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">. Facebook sida . </w:t></w:r>' +
    '</w:p>'
$newPara.Range.InsertXML($xml)

# ---------------------------------------------------------------
# 5) "Beskrivning av webbplatsen" - the _GoBack bookmark now sits
#    right after this run (it used to be in the "Målet..." paragraph).
# ---------------------------------------------------------------
$p = Find-ParaByText $d "Beskrivning av webbplatsen"
$xml = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:pPr>' +
    '<w:pStyle w:val="Rubrik2"/>' +
    '<w:widowControl w:val="0"/>' +
    '<w:spacing w:before="360" w:after="80"/>' +
    '<w:ind w:right="210"/>' +
    '<w:jc w:val="both"/>' +
    '</w:pPr>' +
    '<w:bookmarkStart w:id="100" w:name="_6jzzhl6plkw8" w:colFirst="0" w:colLast="0"/>' +
    '<w:bookmarkEnd w:id="100"/>' +
    '<w:r><w:t>Beskrivning av webbplatsen</w:t></w:r>' +
    '<w:bookmarkStart w:id="101" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="101"/>' +
    '</w:p>'
$p.Range.InsertXML($xml)

Write-Host "done"
